# =====================================================================
# Applies the "Commit On Sep 25" changes to AccountCreationDetails.xlsx
#   - Sheet1: replace the sample data with new sample data + hyperlinks
#   - Adds Sheet2 and Sheet3 with new sample data + hyperlinks
#   - Updates selections / active tab to match the final authored state
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# SHEET1 - update existing data
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# remove the two pre-existing hyperlinks so we can rebuild them (and their
# relationship ids) in the desired final order
$ws1.Hyperlinks.Delete()

# Row 1 (headers)
$ws1.Range("A1").Value = "firstname"
$ws1.Range("B1").Value = "lastname"
$ws1.Range("C1").Value = "phone"
$ws1.Range("D1").Value = "email"
$ws1.Range("E1").Value = "password"

# Row 2
$ws1.Range("A2").Value = "abcd"
$ws1.Range("B2").Value = 1234
$ws1.Range("C2").Value = "564.625uu"
$ws1.Range("D2").Value = "noone@hotmail.com"
$ws1.Range("E2").Value = "abcd12343"

# Row 3
$ws1.Range("A3").Value = "Jinesh"
$ws1.Range("B3").Value = "Z@lawadia"
$ws1.Range("C3").Value = "456789ii"
$ws1.Range("D3").Value = "email@email.com"
$ws1.Range("E3").Value = "459766…"

# Row 4
$ws1.Range("A4").Value = "Pari"
$ws1.Range("B4").Value = "Zalawadia"
$ws1.Range("C4").Value = 6135698742
$ws1.Range("D4").Value = "pari.patel4192@gmail.com"
$ws1.Range("E4").Value = "abcd12343"

# Hyperlinks (added in order so relationship ids come out rId1..rId4)
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:noone@hotmail.com") | Out-Null
$ws1.Range("D2").Style = "Hyperlink"

$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:Z@lawadia") | Out-Null
$ws1.Range("B3").Style = "Hyperlink"

$ws1.Hyperlinks.Add($ws1.Range("D3"), "mailto:email@email.com") | Out-Null
$ws1.Range("D3").Style = "Hyperlink"

$ws1.Hyperlinks.Add($ws1.Range("D4"), "mailto:pari.patel4192@gmail.com") | Out-Null
$ws1.Range("D4").Style = "Hyperlink"

# Column widths (closest representable values to the authored widths)
$ws1.Columns.Item(3).ColumnWidth = 10.083333333352488   # -> 11,    bestFit
$ws1.Columns.Item(4).ColumnWidth = 16.583333333340946   # -> 17.44..
$ws1.Columns.Item(5).ColumnWidth = 10.083333333352488   # -> 11

$ws1.Range("D1").Select() | Out-Null

# ---------------------------------------------------------------
# SHEET2 - brand new sheet, placed right after Sheet1
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Range("A1").Value = "firstname"
$ws2.Range("B1").Value = "lastname"
$ws2.Range("C1").Value = "postalcode"
$ws2.Range("D1").Value = "email"
$ws2.Range("E1").Value = "password"
$ws2.Range("F1").Value = "current password"
$ws2.Range("G1").Value = "new password"

$ws2.Range("A2").Value = "Pari"
$ws2.Range("B2").Value = "Zalawadia"
$ws2.Range("C2").Value = "M6L 1B4"
$ws2.Range("D2").Value = "abcd@gmail.com"
$ws2.Range("E2").Value = "Jinesh@2694"
$ws2.Range("F2").Value = "Jinesh@2694"
$ws2.Range("G2").Value = "abcd_1234"

$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:abcd@gmail.com") | Out-Null
$ws2.Range("D2").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:Jinesh@2694") | Out-Null
$ws2.Range("E2").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("F2"), "mailto:Jinesh@2694") | Out-Null
$ws2.Range("F2").Style = "Hyperlink"

$ws2.Columns.Item(3).ColumnWidth = 12.750000000009788   # -> 13.6640625
$ws2.Columns.Item(4).ColumnWidth = 17.250000000020325   # -> 18.109375
$ws2.Columns.Item(5).ColumnWidth = 15.083333333357508   # -> 16
$ws2.Columns.Item(6).ColumnWidth = 20.916666666696695   # -> 21.88671875
$ws2.Columns.Item(7).ColumnWidth = 11.25000000002635    # -> 12.21875
$ws2.Columns.Item(8).ColumnWidth = 16.583333333340946   # -> 17.5546875

$ws2.Range("H5").Select() | Out-Null

# ---------------------------------------------------------------
# SHEET3 - brand new sheet, placed right after Sheet2 (becomes active)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)

$ws3.Range("A1").Value = "first name"
$ws3.Range("B1").Value = "last name"
$ws3.Range("C1").Value = "address line 1"
$ws3.Range("D1").Value = "city name"
$ws3.Range("E1").Value = "postal code"
$ws3.Range("F1").Value = "phone number"

$ws3.Range("A2").Value = "Pari"
$ws3.Range("B2").Value = "Zalawadia"
$ws3.Range("C2").Value = "59 Bayshore drive"
$ws3.Range("D2").Value = "Ottawa"
$ws3.Range("E2").Value = "M2L 1L5"
$ws3.Range("F2").Value = 613456789

$ws3.Columns.Item(3).ColumnWidth = 15.916666666691677   # -> 16.77734375
$ws3.Columns.Item(6).ColumnWidth = 9.083333333334641    # -> 10, bestFit

$ws3.Range("F3").Select() | Out-Null

Write-Host "edit complete"
